# Correccion de subida de archivos
# Remove the "IdServiceTypes" column (P) which held the "Terrestre nacional"
# constant, shifting the "status" column (Q) left into its place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(16).Delete()

# The saved file also shows the "status" column (now P) selected in its
# entirety, matching what a user would see right after deleting column P.
$ws.Range("P1:P1048576").Select()
